# Updated cryptos list on Tue Oct 10 09:10:42 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns with the latest
# scraped figures, and corrects the Chainlink / ShibaInu row ordering
# (rows 19-20 had swapped places).
#
# Values are assigned with a leading apostrophe to force literal-text
# interpretation (several prices, e.g. "209.20" or "7.38", would otherwise
# be auto-converted to numbers by Excel's type inference), then the cell
# Style is reset to "Normal" so no stray text-format style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.761.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.27%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.595.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.55%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.15%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'209.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.97%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -1.82%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.16%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'22.36"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.49%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.31%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.60%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.53%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.822.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.47%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.574.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.40%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.50%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.532"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.49%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'27.763.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'63.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.55%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'219.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.13%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'Chainlink"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'7.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.54%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'ShibaInu"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0696"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.07%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.15%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -3.24%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'9.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.74%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -4.22%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'153.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.19%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.83%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.14%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.86%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.105"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.45%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.69%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.09%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -4.08%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.376.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.54%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -2.92%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.56%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.978"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.38%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.05%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0169"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.18%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -2.64%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.828"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.67%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.03%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.973"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.62%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'64.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.84%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +2.67%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'5.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.03%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.69%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.733.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.47%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'86.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.07%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.61%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0966"
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.97%  "
$ws.Range("E51").Style = "Normal"
